# Append three new log rows (2018.07.04 measurements) to the bottom of the
# log table, matching the new data referenced by the updated paper.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that looks like a date/number (e.g. "2018.07.04")
# into a cell while forcing it to stay as plain text, exactly like the
# other "Date" column entries already in the sheet, and without leaving
# a non-default style on the cell.
function Set-TextValue($cell, [string]$text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$newRows = @(
    @{ Date = "2018.07.04"; Time = "15:06:46"; Neuron = "RS"; Astim = 10; Tstim = 250; PRF = 0.1;   Duty = 0.98; Samples = 3475; CompTime = 0.4; Spikes = 0; Latency = "N/A"; Rate = "N/A" },
    @{ Date = "2018.07.04"; Time = "15:06:52"; Neuron = "RS"; Astim = 10; Tstim = 250; PRF = 0.107; Duty = 0.98; Samples = 3511; CompTime = 0.4; Spikes = 0; Latency = "N/A"; Rate = "N/A" },
    @{ Date = "2018.07.04"; Time = "15:07:02"; Neuron = "RS"; Astim = 6;  Tstim = 250; PRF = 0.107; Duty = 0.98; Samples = 3511; CompTime = 0.3; Spikes = 0; Latency = "N/A"; Rate = "N/A" }
)

$startRow = 24
$r = $startRow
foreach ($row in $newRows) {
    Set-TextValue $ws.Cells.Item($r, 1) $row.Date

    $ws.Cells.Item($r, 2).Value = $row.Time
    $ws.Cells.Item($r, 3).Value = $row.Neuron
    $ws.Cells.Item($r, 4).Value = $row.Astim
    $ws.Cells.Item($r, 5).Value = $row.Tstim
    $ws.Cells.Item($r, 6).Value = $row.PRF
    $ws.Cells.Item($r, 7).Value = $row.Duty
    $ws.Cells.Item($r, 8).Value = $row.Samples
    $ws.Cells.Item($r, 9).Value = $row.CompTime
    $ws.Cells.Item($r, 10).Value = $row.Spikes
    $ws.Cells.Item($r, 11).Value = $row.Latency
    $ws.Cells.Item($r, 12).Value = $row.Rate

    $r++
}
